# Trade #41 closed at 2026-02-17 15:28:57 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet -----------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.25
$summary.Range("B4").Value = 0.25
$summary.Range("B6").Value = 41
$summary.Range("B7").Value = 13
$summary.Range("B9").Value = 31.71

# --- Strategy Status sheet ---------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.25
$status.Range("D4").Value = 41
$status.Range("E4").Value = 0.25
$status.Range("F4").Value = 0.25
$status.Range("G4").Value = 31.71

# --- New trade row shared by "All Trades" and "MarketMaking" sheets ----
$tradeRow = @{
    A = 41
    B = "2026-02-17"
    C = "15:28:51"
    D = "MarketMaking"
    E = "UP"
    F = 0.02
    G = 0.03
    H = "CLOSED"
    I = 50
    J = 0.01
    K = 100.25
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 42
    $ws.Cells.Item($newRow, 1).Value = $tradeRow.A

    # Column B holds a date-formatted string ("2026-02-17"). A plain .Value
    # assignment gets auto-parsed into a real Excel date serial, which does
    # not match the source data (stored as literal text). Pre-formatting the
    # cell as Text keeps the literal string, then resetting the cell style
    # back to Normal avoids leaving a stray "@" number format behind.
    $bCell = $ws.Cells.Item($newRow, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $tradeRow.B
    $bCell.Style = "Normal"

    $ws.Cells.Item($newRow, 3).Value = $tradeRow.C
    $ws.Cells.Item($newRow, 4).Value = $tradeRow.D
    $ws.Cells.Item($newRow, 5).Value = $tradeRow.E
    $ws.Cells.Item($newRow, 6).Value = $tradeRow.F
    $ws.Cells.Item($newRow, 7).Value = $tradeRow.G
    $ws.Cells.Item($newRow, 8).Value = $tradeRow.H
    $ws.Cells.Item($newRow, 9).Value = $tradeRow.I
    $ws.Cells.Item($newRow, 10).Value = $tradeRow.J
    $ws.Cells.Item($newRow, 11).Value = $tradeRow.K
    $ws.Cells.Item($newRow, 12).Value = $tradeRow.L
    $ws.Cells.Item($newRow, 13).Value = $tradeRow.M
    $ws.Cells.Item($newRow, 14).Value = $tradeRow.N
    $ws.Cells.Item($newRow, 15).Value = $tradeRow.O
    $ws.Cells.Item($newRow, 16).Value = $tradeRow.P
    $ws.Cells.Item($newRow, 17).Value = $tradeRow.Q
}
